# Weekly price-sheet update: a new record for "Feria Lagunitas de Puerto Montt -
# Pepino ensalada" is inserted at row 219, pushing the existing rows 219-342
# down to 220-343 (the sheet's used range grows from A1:R342 to A1:R343).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 219; this shifts rows 219:342 down to 220:343 and
# carries their formatting (including the date number-format on column D).
$ws.Rows(219).Insert()

# Populate the newly inserted row 219 with this week's record. The
# categorical / descriptive columns are identical to every other row for
# this market+product combination; only the date (D), volume (J), min/max/avg
# price (K/L/M) and $/Kg (P) are new for this entry.
$ws.Range("A219").Value = 4
$ws.Range("B219").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C219").Value = "Los Lagos"
$ws.Range("D219").Value = 44879
$ws.Range("E219").Value = 10
$ws.Range("F219").Value = 100112043
$ws.Range("G219").Value = "Pepino ensalada"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 200
$ws.Range("K219").Value = 26000
$ws.Range("L219").Value = 26000
$ws.Range("M219").Value = 26000
$ws.Range("N219").Value = "$/caja 60 unidades"
$ws.Range("O219").Value = "Región de Arica y Parinacota"
$ws.Range("P219").Value = 433
$ws.Range("Q219").Value = 60
$ws.Range("R219").Value = "Hortaliza"
